# Swap the presentation's theme colour palette from the "Integral" /
# "Red Violet" scheme to the "Office Theme" / "Office" scheme.
#
# The PowerPoint object model exposes the 12 DrawingML theme colours
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) through
# Slide.ThemeColorScheme.Colors(1..12) as COM RGB (BGR-packed) integers.
# Writing ThemeColor.RGB updates the presentation's <a:clrScheme> in
# ppt/theme/theme1.xml (the theme used by the slide master / all slides).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette ("Office Theme"), in clrScheme order 1-12.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

function ConvertTo-ComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

for ($i = 1; $i -le 12; $i++) {
    $color = $tcs.Colors($i)
    $color.RGB = ConvertTo-ComRgb($officeColors[$i - 1])
}
